# Overall_Rebate_Efficiency.xlsx update:
# Append 5 new weekly rows (Week_31..Week_35) to the OverallRebateEfficiency sheet
# and restore the first sheet as the active/selected sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("OverallRebateEfficiency")

# New weekly rows to append after the existing data (rows 2-30).
$newRows = @(
    @{ Week = "Week_31"; Value = 0.734 },
    @{ Week = "Week_32"; Value = 0.862 },
    @{ Week = "Week_33"; Value = 0.804 },
    @{ Week = "Week_34"; Value = 0.7631 },
    @{ Week = "Week_35"; Value = 0.7837 }
)

$startRow = 31
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws1.Range("A" + $r).Value = $newRows[$i].Week
    $ws1.Range("B" + $r).Value = $newRows[$i].Value
}

# Make the OverallRebateEfficiency sheet the active sheet/tab and update the
# visible selection to match the new bottom of the data range.
$ws1.Activate()
$ws1.Range("C35").Select()
